$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" / "is_enabled" comment columns (F and G), shifting
# the trailing "rem" column left into F.
$ws.Range("F1:G1").EntireColumn.Delete()
